# Initial Time.xlsx update (WRI China -> Hong Kong EPS v2.0.0 bring-up)
# Core content change: the "Year" value on the IT sheet advances from 2016 to 2018.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsIT    = $wb.Worksheets.Item("IT")

# Update the Initial Time year.
$wsIT.Range("B2").Value = 2018

# Leave the IT sheet with cell B3 highlighted, then switch back to the
# "About" sheet (which becomes the active/selected tab) with A21 selected,
# matching the saved view state in the workbook.
$wsIT.Range("B3").Select() | Out-Null

$wsAbout.Activate() | Out-Null
$wsAbout.Range("A21").Select() | Out-Null
